$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("103:103").Insert()

$ws.Range("A103").Value = 4
$ws.Range("B103").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C103").Value = "Los Lagos"
$ws.Range("D103").Value = 44694
$ws.Range("E103").Value = 10
$ws.Range("F103").Value = 100112052
$ws.Range("G103").Value = "Albahaca"
$ws.Range("H103").Value = "Sin especificar"
$ws.Range("I103").Value = "Primera"
$ws.Range("J103").Value = 90
$ws.Range("K103").Value = 8000
$ws.Range("L103").Value = 8000
$ws.Range("M103").Value = 8000
$ws.Range("N103").Value = "$/docena de matas"
$ws.Range("O103").Value = "Región Metropolitana"
$ws.Range("P103").Value = 1333
$ws.Range("Q103").Value = 6
$ws.Range("R103").Value = "Hortaliza"
